$wb = $excel.ActiveWorkbook

# --- 1. PUC sheet: add new account code row (9605 / "Integracion diagonal ") ---
$ws1 = $wb.Worksheets.Item("PUC")
$ws1.Range("B2476").Value = 9605
$ws1.Range("C2476").Value = "Integracion diagonal "

# --- 2. cuentas_modelo sheet / Tabla1: append the matching row ---
$ws2 = $wb.Worksheets.Item("cuentas_modelo")
$lo = $ws2.ListObjects.Item("Tabla1")

$ws2.Range("A16").Value = 9605
$lo.Resize($ws2.Range("A1:E16"))

$ws2.Range("B16").Formula = '=_xlfn.CONCAT(VALUE(LEFT(A16,1))," ",VLOOKUP(VALUE(LEFT(A16,1)),PUC!$B$3:$C$2475,2,0))'
$ws2.Range("C16").Formula = '=+VALUE(LEFT(A16,1))'
$ws2.Range("D16").Formula = '=_xlfn.CONCAT(VALUE(LEFT(A16,2))," ",VLOOKUP(VALUE(LEFT(A16,2)),PUC!$B$3:$C$2475,2,0))'
$ws2.Range("E16").Formula = '=+VLOOKUP(Tabla1[[#This Row],[codigo_cuenta]],PUC!B17:C2489,2,0)'

# --- 3. Sort Tabla1 ascending by codigo_cuenta (was descending) ---
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws2.Range("A2:A16"))
$lo.Sort.Header = 1
$lo.Sort.Apply()

# --- 4. Window / selection state: sheet2 becomes the active tab ---
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 2465
$ws1.Range("C2477").Select()

$ws2.Activate()
$ws2.Range("A16").Select()
